$d = $word.ActiveDocument

$replacements = @(
    @("733×6=", "505×2="),
    @("926×8=", "538×9="),
    @("556×9=", "597×8="),
    @("635×4=", "604×8="),
    @("381×8=", "483×6="),
    @("853×6=", "261×4="),
    @("808×4=", "803×6="),
    @("179×7=", "600×7="),
    @("613×8=", "841×5="),
    @("999×9=", "714×9="),
    @("980×9=", "446×7="),
    @("561×4=", "763×2="),
    @("753×2=", "952×8="),
    @("257×3=", "929×7="),
    @("643×7=", "272×7="),
    @("632×4=", "892×8="),
    @("711×8=", "748×4="),
    @("902×9=", "804×3="),
    @("352×5=", "460×7="),
    @("465×4=", "248×9="),
    @("544×6=", "933×4="),
    @("680×9=", "726×5="),
    @("307×7=", "630×3="),
    @("516×3=", "561×2="),
    @("602×3=", "268×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
